$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cells with new values
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 5

$ws.Range("A4").Value = 4
$ws.Range("D4").Value = 7

$ws.Range("D5").Value = 5

$ws.Range("A6").Value = 5
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 10

$ws.Range("A7").Value = 7
$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 15

$ws.Range("A8").Value = 3
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5

# Add new row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 10
